$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row 20 with value 12345 in column A
$ws.Range("A20").Value = 12345

# Set the active cell / selection to A2 (as reflected in the saved view state)
$ws.Range("A2").Select()
